$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column formatting (style index of A238) down through the new rows
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)

$dates = @(44313, 44314, 44315, 44316, 44317, 44318)
$bvals = @(0, 0, 0, 0, 1, 0)
$cvals = @(0, 0, 0, 0, 1, 1)
$dvals = @(0, 0, 0, 0, 83.40283569641367, 83.40283569641367)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 239 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    $ws.Cells.Item($row, 3).Value = $cvals[$i]
    $ws.Cells.Item($row, 4).Value = $dvals[$i]
}
